$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A70").Value = "Kamel"
$ws.Range("B70").Value = "Mhalhel"
$ws.Range("C70").Value = "Università Degli Studi di Messina"
$ws.Range("D70").Value = "Italie"
$ws.Range("E70").Value = "R35oeVAAAAAJ"
$ws.Range("F70").Value = "M"
$ws.Range("G70").Value = 1992
$ws.Range("H70").Value = "Médecine, Biologie et Sciences de la Santé"

$ws.Range("F69").Copy()
$ws.Range("F70").PasteSpecial(-4122)

$ws.Range("H71").Select()
